# Weekly update: insert a new price-report row for Kiwi / Hayward / Especial
# (Región de O'Higgins, $/bandeja 10 kilos) ahead of the existing row 310,
# pushing the rest of the data set (rows 310-360) down to 311-361.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 310 - everything below shifts down by one.
$ws.Rows.Item(310).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A310").Value = 10
$ws.Range("B310").Value = "Vega Modelo de Temuco"
$ws.Range("C310").Value = "La Araucanía"
$ws.Range("D310").Value = 44504
$ws.Range("E310").Value = 9
$ws.Range("F310").Value = "Fruta"
$ws.Range("G310").Value = 100101
$ws.Range("H310").Value = "Berries"
$ws.Range("I310").Value = 100101007
$ws.Range("J310").Value = "Kiwi"
$ws.Range("K310").Value = "Hayward"
$ws.Range("L310").Value = "Especial"
$ws.Range("M310").Value = 200
$ws.Range("N310").Value = 16000
$ws.Range("O310").Value = 16000
$ws.Range("P310").Value = 16000
$ws.Range("Q310").Value = "$/bandeja 10 kilos"
$ws.Range("R310").Value = "Región de O'Higgins"
$ws.Range("S310").Value = 1600
$ws.Range("T310").Value = 10
